$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 114.5
$ws.Range("I4").Value = 73.5
$ws.Range("K4").Value = 73.5
$ws.Range("M4").Value = 40.5
# Row 29
$ws.Range("H29").Value = 7945.7827
$ws.Range("J29").Value = 8997.375
$ws.Range("L29").Value = 26992.125
$ws.Range("N29").Value = -27554.125
# Row 58
$ws.Range("H58").Value = 10281.8
$ws.Range("I58").Value = 10281.8
$ws.Range("K58").Value = 30845.4
$ws.Range("M58").Value = -30695.4
# Row 69
$ws.Range("H69").Value = 4500
$ws.Range("I69").Value = 4500
$ws.Range("K69").Value = 13500
$ws.Range("M69").Value = -12626
# Row 72
$ws.Range("H72").Value = 4500
$ws.Range("I72").Value = 4500
$ws.Range("K72").Value = 40500
$ws.Range("M72").Value = -36132
# Row 80
$ws.Range("H80").Value = 1537.9615
$ws.Range("I80").Value = 1234.5294
$ws.Range("J80").Value = 2111.111
$ws.Range("K80").Value = 3703.5882
$ws.Range("L80").Value = 6333.333
$ws.Range("M80").Value = -2705.5882
$ws.Range("N80").Value = -8329.332999999999
# Row 83
$ws.Range("H83").Value = 1537.9615
$ws.Range("I83").Value = 1234.5294
$ws.Range("J83").Value = 2111.111
$ws.Range("K83").Value = 11110.7646
$ws.Range("L83").Value = 18999.999
$ws.Range("M83").Value = -6118.764599999999
$ws.Range("N83").Value = -28983.999
# Row 106
$ws.Range("H106").Value = 25322.615
$ws.Range("I106").Value = 29489.818
$ws.Range("K106").Value = 29489.818
$ws.Range("M106").Value = -28858.818
# Row 111
$ws.Range("H111").Value = 4043.125
$ws.Range("I111").Value = 3881.25
$ws.Range("K111").Value = 11643.75
$ws.Range("M111").Value = -8576.75
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("N120").ClearContents()
# Row 132
$ws.Range("H132").Value = 1371.2963
$ws.Range("I132").Value = 1308.7307
$ws.Range("K132").Value = 3926.1921
$ws.Range("M132").Value = -1396.1921
# Row 137
$ws.Range("H137").Value = 9263962
$ws.Range("I137").Value = 14706944
$ws.Range("J137").Value = 10892.5
$ws.Range("K137").Value = 44120832
$ws.Range("L137").Value = 32677.5
$ws.Range("M137").Value = -44118282
$ws.Range("N137").Value = -37777.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 29022.918
$ws.Range("I32").Value = 29812.8
$ws.Range("K32").Value = 29812.8
$ws.Range("M32").Value = -29525.8
# Row 50
$ws.Range("H50").Value = 1983.6666
$ws.Range("J50").Value = 1161.3334
$ws.Range("L50").Value = 1161.3334
$ws.Range("N50").Value = -2589.3334
# Row 110
$ws.Range("H110").Value = 14707795
$ws.Range("J110").Value = 3372.5
$ws.Range("L110").Value = 3372.5
$ws.Range("N110").Value = -7462.5
# Row 132
$ws.Range("H132").Value = 6663.5454
$ws.Range("I132").Value = 4533.3887
$ws.Range("K132").Value = 13600.1661
$ws.Range("M132").Value = -11070.1661

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 283668
$ws.Range("I86").Value = 2971.9092
$ws.Range("K86").Value = 2971.9092
$ws.Range("M86").Value = -1848.9092
# Row 89
$ws.Range("H89").Value = 283668
$ws.Range("I89").Value = 2971.9092
$ws.Range("K89").Value = 14859.546
$ws.Range("M89").Value = -9243.546
# Row 134
$ws.Range("H134").Value = 5929.7407
$ws.Range("I134").Value = 2567.7334
$ws.Range("K134").Value = 7703.2002
$ws.Range("M134").Value = -5168.2002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 6129.8887
$ws.Range("I16").Value = 5254.4
$ws.Range("J16").Value = 7224.25
$ws.Range("K16").Value = 5254.4
$ws.Range("L16").Value = 7224.25
$ws.Range("M16").Value = -4967.4
$ws.Range("N16").Value = -7798.25
# Row 105
$ws.Range("H105").Value = 2151.125
$ws.Range("I105").Value = 2151.125
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2151.125
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -404.125
$ws.Range("N105").ClearContents()
# Row 107
$ws.Range("H107").Value = 1299.2
$ws.Range("I107").Value = 1499
$ws.Range("K107").Value = 1499
$ws.Range("M107").Value = 421
# Row 113
$ws.Range("H113").Value = 6129.8887
$ws.Range("I113").Value = 5254.4
$ws.Range("J113").Value = 7224.25
$ws.Range("K113").Value = 5254.4
$ws.Range("L113").Value = 7224.25
$ws.Range("M113").Value = -3084.4
$ws.Range("N113").Value = -11564.25
# Row 132
$ws.Range("H132").Value = 38853.406
$ws.Range("I132").Value = 4747.6333
$ws.Range("K132").Value = 14242.8999
$ws.Range("M132").Value = -11712.8999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 460.5
$ws.Range("I17").Value = 460.5
$ws.Range("K17").Value = 1381.5
$ws.Range("M17").Value = -1212.5
# Row 81
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 30000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -32246
# Row 84
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 90000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -101232
# Row 119
$ws.Range("H119").Value = 402.4
$ws.Range("I119").Value = 402.4
$ws.Range("K119").Value = 1207.2
$ws.Range("M119").Value = 3630.8
# Row 122
$ws.Range("H122").Value = 10741
$ws.Range("J122").Value = 989.5
$ws.Range("L122").Value = 8905.5
$ws.Range("N122").Value = -13805.5
# Row 131
$ws.Range("H131").Value = 17552174
$ws.Range("I131").Value = 83334470
$ws.Range("J131").Value = 10227.934
$ws.Range("K131").Value = 250003410
$ws.Range("L131").Value = 30683.802
$ws.Range("M131").Value = -249998370
$ws.Range("N131").Value = -40763.802

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 11285.75
$ws.Range("I122").Value = 11285.75
$ws.Range("K122").Value = 33857.25
$ws.Range("M122").Value = -31407.25
# Row 126
$ws.Range("H126").Value = 2973.8438
$ws.Range("I126").Value = 2270.875
$ws.Range("J126").Value = 3676.8125
$ws.Range("K126").Value = 6812.625
$ws.Range("L126").Value = 11030.4375
$ws.Range("M126").Value = -4342.625
$ws.Range("N126").Value = -15970.4375
# Row 132
$ws.Range("H132").Value = 10357.523
$ws.Range("I132").Value = 10049.3125
$ws.Range("K132").Value = 30147.9375
$ws.Range("M132").Value = -27617.9375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3183.1292
$ws.Range("I22").Value = 1950.6875
$ws.Range("J22").Value = 4497.7334
$ws.Range("K22").Value = 1950.6875
$ws.Range("L22").Value = 4497.7334
$ws.Range("M22").Value = -1655.6875
$ws.Range("N22").Value = -5087.7334
# Row 27
$ws.Range("H27").Value = 3183.1292
$ws.Range("I27").Value = 1950.6875
$ws.Range("J27").Value = 4497.7334
$ws.Range("K27").Value = 1950.6875
$ws.Range("L27").Value = 4497.7334
$ws.Range("M27").Value = -1843.6875
$ws.Range("N27").Value = -4711.7334
# Row 100
$ws.Range("J100").Value = 2586.2727
$ws.Range("L100").Value = 2586.2727
$ws.Range("N100").Value = -3668.2727
# Row 132
$ws.Range("H132").Value = 4636.8125
$ws.Range("I132").Value = 2942.1428
$ws.Range("J132").Value = 16499.5
$ws.Range("K132").Value = 8826.428400000001
$ws.Range("L132").Value = 49498.5
$ws.Range("M132").Value = -6296.428400000001
$ws.Range("N132").Value = -54558.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 47497.445
$ws.Range("J54").Value = 28434.625
$ws.Range("L54").Value = 28434.625
$ws.Range("N54").Value = -29474.625
# Row 62
$ws.Range("H62").Value = 7983.3335
$ws.Range("I62").Value = 7983.3335
$ws.Range("K62").Value = 7983.3335
$ws.Range("M62").Value = -7359.3335
# Row 65
$ws.Range("H65").Value = 7983.3335
$ws.Range("I65").Value = 7983.3335
$ws.Range("K65").Value = 39916.6675
$ws.Range("M65").Value = -36796.6675
# Row 81
$ws.Range("H81").Value = 6016.478
$ws.Range("I81").Value = 1807.3846
$ws.Range("K81").Value = 3614.7692
$ws.Range("M81").Value = -2553.7692
# Row 84
$ws.Range("H84").Value = 6016.478
$ws.Range("I84").Value = 1807.3846
$ws.Range("K84").Value = 18073.846
$ws.Range("M84").Value = -12769.846
# Row 126
$ws.Range("H126").Value = 4812.6665
$ws.Range("I126").Value = 3199.9092
$ws.Range("J126").Value = 9247.75
$ws.Range("K126").Value = 9599.7276
$ws.Range("L126").Value = 27743.25
$ws.Range("M126").Value = -7129.7276
$ws.Range("N126").Value = -32683.25
# Row 132
$ws.Range("H132").Value = 4825.8374
$ws.Range("I132").Value = 2497.0344
$ws.Range("K132").Value = 7491.1032
$ws.Range("M132").Value = -4961.1032
